# Apply the "contingencies with rene fine" edit:
#  - Extend the data table from columns B:O to B:Q (two new columns, 14 and 15)
#  - Swap the values in columns I/K and M/O for every data row (2-25)
#  - Populate the two new columns (P, Q) with the value 2 for every data row (2-25)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 = 14, Q1 = 15, matching the formatting of O1 ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows (2-25): swap I<->K values and M<->O values ---
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# --- Data rows (2-25): new columns P and Q, all set to 2 ---
$ws.Range("P2:Q25").Value = 2
